# Largest Range - Solution 1 (Can be done better)
# Adds two new rows to the AlgoExpert tracking sheet:
#   Row 5: "Vallidate BST" problem
#   Row 6: "Largrest Range" problem
# plus two new review comments (E5, F5) and assorted wrap-text tweaks that
# Excel picked up on existing cells while the author was editing the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# So that newly authored comments attribute to the sheet's usual commenter.
$excel.UserName = "Gokul Chagalamarri Nippani"

# --- Wrap-text touch-ups on pre-existing cells -----------------------------
$ws.Range("B1").WrapText = $true
$ws.Range("C1").WrapText = $true
$ws.Range("B2").WrapText = $true
$ws.Range("D3").WrapText = $true
$ws.Range("B4").WrapText = $true

# --- Row 5: Vallidate BST ---------------------------------------------------
$ws.Range("A5").Value = "Vallidate BST"

$ws.Range("B5").Value = "https://www.geeksforgeeks.org/a-program-to-check-if-a-binary-tree-is-bst-or-not/"
$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.geeksforgeeks.org/a-program-to-check-if-a-binary-tree-is-bst-or-not/") | Out-Null
$ws.Range("B5").WrapText = $true

$ws.Range("C5").Value = "We need to compare each node value with max and min value "
$ws.Range("C5").WrapText = $true

$ws.Range("D5").Value = "initialize max and min value"
$ws.Range("D5").Font.Bold = $false

$ws.Range("E5").Value = "O(n)"
$ws.Range("E5").Style = "Good"

$ws.Range("F5").Value = "O(d)"

$ws.Rows(5).RowHeight = 28.8

# --- Row 6: Largrest Range ---------------------------------------------------
$ws.Range("A6").Value = "Largrest Range"

$ws.Range("C6").Value = "Solution 1: Make sure to compare current range with largest range before returning the result."
$ws.Range("C6").WrapText = $true

$ws.Range("E6").Value = "O(nlogn)"
$ws.Range("E6").Style = "Bad"

$ws.Range("F6").Value = "O(1) "
$ws.Range("F6").Style = "Good"

$ws.Range("G6").Value = "O(n)"
$ws.Range("G6").Style = "Good"

$ws.Range("H6").Value = "O(n)"

$ws.Rows(6).RowHeight = 43.2

# --- New review comments ----------------------------------------------------
$ws.Range("E5").AddComment("Gokul Chagalamarri Nippani:`nWe will be accessing all the nodes.") | Out-Null
$ws.Range("F5").AddComment("Gokul Chagalamarri Nippani:`nCall Stack. Where d is depth of tree.") | Out-Null

# --- Selection as left by the author after the edit -------------------------
$ws.Range("D12").Select() | Out-Null
